$d = $word.ActiveDocument

# Full package XML (includes all parts: document.xml, numbering.xml, etc.)
$full = $d.Content.WordOpenXML

# ---------------------------------------------------------------------------
# 1) document.xml: insert the new "Draft 3" section right before the final
#    empty paragraph (w14:paraId="6E156D67").
# ---------------------------------------------------------------------------
$newParagraphs = @'
<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></w:rPr><w:lastRenderedPageBreak/><w:t>✅</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:b/><w:bCs/></w:rPr><w:t>Draft 3: Design-Oriented Emphasis (Ideal for visual-focused documentation)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>4.1 User Interfaces</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Syarti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> system features a visually engaging, RTL-based single-page web interface designed for customers renting vehicles. The prototype showcases a clean and intuitive design built with Bootstrap 5 RTL, HTML5, CSS3, and Font Awesome.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/></w:pPr><w:r><w:t>The layout includes:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">A dark-themed </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>sticky navigation bar</w:t></w:r><w:r><w:t xml:space="preserve"> with icon-enhanced links.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">A promotional </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>hero section</w:t></w:r><w:r><w:t xml:space="preserve"> with a prominent booking call-to-action.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">An interactive </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>search card</w:t></w:r><w:r><w:t xml:space="preserve"> allowing customers to choose pickup locations and rental dates.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>Vehicle cards</w:t></w:r><w:r><w:t xml:space="preserve"> displaying cars with labeled badges (e.g., “</w:t></w:r><w:r><w:rPr><w:rtl/></w:rPr><w:t>عرض خاص</w:t></w:r><w:r><w:t>”) and specifications.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>testimonial section</w:t></w:r><w:r><w:t xml:space="preserve"> with avatar images, quotes, and rating icons.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="23"/></w:numPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/></w:rPr><w:t>footer</w:t></w:r><w:r><w:t xml:space="preserve"> containing social links and brief company information.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/></w:pPr><w:r><w:t>All buttons, including “</w:t></w:r><w:r><w:rPr><w:rtl/></w:rPr><w:t>احجز الآن</w:t></w:r><w:r><w:t>” and navigation links, point to placeholder anchors. Accessibility features have not yet been integrated. The UI has been tested on smartphones and desktops for responsiveness.</w:t></w:r></w:p>
'@

$anchorDoc = '<w:p w14:paraId="6E156D67"'
if ($full.IndexOf($anchorDoc) -lt 0) {
    throw "document.xml anchor not found"
}
$full = $full.Replace($anchorDoc, $newParagraphs + $anchorDoc)

# ---------------------------------------------------------------------------
# 2) numbering.xml: insert a new abstractNum (nsid 15B3196F) as id "1",
#    renumber the existing abstractNum ids 1->2 and 2->3, and add/renumber
#    the <w:num> entries so numId 23 is the new bulleted list.
# ---------------------------------------------------------------------------
$newAbstractNum = @'
<w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="15B3196F"/><w:multiLevelType w:val="multilevel"/><w:tmpl w:val="F990A650"/><w:lvl w:ilvl="0"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="720"/></w:tabs><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="1440"/></w:tabs><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="2160"/></w:tabs><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="2880"/></w:tabs><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="3600"/></w:tabs><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="4320"/></w:tabs><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="5040"/></w:tabs><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="5760"/></w:tabs><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="6480"/></w:tabs><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/><w:sz w:val="20"/></w:rPr></w:lvl></w:abstractNum>
'@

# 2a. Renumber abstractNum 2 (nsid 5A0132B6) -> 3
$old1 = '<w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="5A0132B6"'
$new1 = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="5A0132B6"'
if ($full.IndexOf($old1) -lt 0) {
    throw "abstractNum 2 (5A0132B6) anchor not found"
}
$full = $full.Replace($old1, $new1)

# 2b. Insert the new abstractNum before abstractNum 1 (nsid 36052345) and renumber it -> 2
$old2 = '<w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="36052345"'
$new2 = $newAbstractNum + '<w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="36052345"'
if ($full.IndexOf($old2) -lt 0) {
    throw "abstractNum 1 (36052345) anchor not found"
}
$full = $full.Replace($old2, $new2)

# 2c. num 21 now points at abstractNum 3; add new num 22 pointing at abstractNum 2
$old3 = '<w:num w:numId="21"><w:abstractNumId w:val="2"/></w:num>'
$new3 = '<w:num w:numId="21"><w:abstractNumId w:val="3"/></w:num><w:num w:numId="22"><w:abstractNumId w:val="2"/></w:num>'
if ($full.IndexOf($old3) -lt 0) {
    throw "num 21 anchor not found"
}
$full = $full.Replace($old3, $new3)

# 2d. old num 22 (abstractNumId 1) becomes num 23 (still abstractNumId 1, now the new list)
$old4 = '<w:num w:numId="22"><w:abstractNumId w:val="1"/></w:num>'
$new4 = '<w:num w:numId="23"><w:abstractNumId w:val="1"/></w:num>'
if ($full.IndexOf($old4) -lt 0) {
    throw "num 22 anchor not found"
}
$full = $full.Replace($old4, $new4)

# ---------------------------------------------------------------------------
# Apply the modified package back onto the document.
# ---------------------------------------------------------------------------
$d.Content.InsertXML($full)

Write-Output "Edit applied successfully"
